$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simulate re-execution of the Katalon RAD test suite: the Date column
# records a fresh run timestamp for every test case that actually executed.
$ws.Range("B2").Value = "Thu Dec 07 21:43:01 EST 2023"
$ws.Range("B3").Value = "Thu Dec 07 21:43:12 EST 2023"
$ws.Range("B5").Value = "Thu Dec 07 21:43:23 EST 2023"

# Remove the "Extension Payments" Tax Type from execution.
$ws.Range("C4").Value = "DONOTRUN"

# Widen column C (it's no longer a tight auto-fit once "DONOTRUN" is there)
# and move the selection to the cell that was just edited.
$ws.Columns.Item(3).ColumnWidth = 15
$ws.Range("C4").Select()
